# The "dict_sheet" worksheet holds key/value pairs written out in insertion
# order: row1 = key1/value1, row2 = key2/value2, row3 = key3/value3.
#
# The commit swaps the first and last rows back to their "restored" order
# (row1 <-> row3), leaving row2 untouched - e.g. a python-3-compatible
# dict ordering fix. Net effect on the visible grid:
#   A1: key1   -> key3
#   B1: value1 -> value3
#   A3: key3   -> key1
#   B3: value3 -> value1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dict_sheet")

$ws.Range("A1").Value = "key3"
$ws.Range("B1").Value = "value3"
$ws.Range("A3").Value = "key1"
$ws.Range("B3").Value = "value1"
